# Rename the "MetaData" tab to "Metadata"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MetaData")
$ws.Name = "Metadata"
$ws.Activate()
